# Azure Costs.xlsx - add a "Week" sheet with a weekly Cosmo/Total breakdown,
# a clustered bar chart, and a reused header banner picture.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Create the new "Week" worksheet after "Itemized" (i.e. as the last tab)
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "Week"

# ---------------------------------------------------------------------------
# 2. Fill in the weekly data
#    Row 1 : week labels (B..L)
#    Row 2 : Cosmo costs
#    Row 3 : Total costs
# ---------------------------------------------------------------------------
$weekLabels = @("13-19", "20-26", "27/1-2/2", "3/2-9/2", "10/3-16/2", "17-2/23-2", "24/2-2-/3", "3/3-9/3", "10/3-16/3", "17/3-23/3", "24/3-30-3")
$cosmoValues = @(3450, 6093, 5923, 5390, 6848, 6643, 5072, 4075, 3023, 4486, 3256)
$totalValues = @(5900, 8629, 8470, 7837, 9284, 9041, 7919, 6759, 5741, 8027, 6793)

$cols = @("B", "C", "D", "E", "F", "G", "H", "I", "J", "K", "L")

$ws.Range("A2").Value = "Cosmo"
$ws.Range("A3").Value = "Total"

# Match the original authoring order: every week column except the second
# ("20-26") first, then go back and fill the skipped one in.
for ($i = 0; $i -lt $cols.Length; $i++) {
    if ($i -eq 1) { continue }
    $ws.Range($cols[$i] + "1").Value = $weekLabels[$i]
}
$ws.Range("C1").Value = $weekLabels[1]

for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "2").Value = $cosmoValues[$i]
    $ws.Range($cols[$i] + "3").Value = $totalValues[$i]
}

# ---------------------------------------------------------------------------
# 3. Add the clustered-column chart (series by row: Cosmo / Total)
# ---------------------------------------------------------------------------
$chartAnchorTL = $ws.Range("A4")
$chartAnchorBR = $ws.Range("G18")
$chartLeft = $chartAnchorTL.Left
$chartTop = $chartAnchorTL.Top
$chartWidth = $chartAnchorBR.Left + $chartAnchorBR.Width - $chartLeft
$chartHeight = $chartAnchorBR.Top + $chartAnchorBR.Height - $chartTop

$co = $ws.ChartObjects().Add($chartLeft, $chartTop, $chartWidth, $chartHeight)
$chart = $co.Chart
$chart.ChartType = 51
$chart.SetSourceData($ws.Range("A1:L3"))
$chart.PlotBy = 1

# ---------------------------------------------------------------------------
# 4. Reuse the existing header banner picture (same as on "Itemized")
# ---------------------------------------------------------------------------
$picAnchorTL = $ws.Range("A19")
$picAnchorBR = $ws.Range("Q28")
$picLeft = $picAnchorTL.Left
$picTop = $picAnchorTL.Top
$picWidth = $picAnchorBR.Left + $picAnchorBR.Width - $picLeft
$picHeight = $picAnchorBR.Top + $picAnchorBR.Height - $picTop

try {
    $ws.Shapes.AddPicture("image1.png", 0, 1, $picLeft, $picTop, $picWidth, $picHeight) | Out-Null
} catch {
    Write-Host "AddPicture failed: $_"
}

Write-Host "Week sheet created with chart and picture"
